# Paul Kamau - Task Sheet
# Insert a repeated header row right after the blank separator row (old
# row 43) so that the second alphabetical block of experts (previously
# starting at row 44) gets its own "Name / Credentials / University /
# Location / image" header, and append the same header block again after
# the final blank separator row at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# --- Insert a new row before the old row 44 (shifts rows 44.. down by 1)
$ws.Rows.Item(44).Insert()

# Copy the formatting of the existing header row (row 1) onto the newly
# inserted row, then fill in the header text values.
$ws.Range("B1:F1").Copy()
$ws.Range("B44:F44").PasteSpecial(-4122)

$ws.Range("B44").Value = "Name"
$ws.Range("C44").Value = "Credentials"
$ws.Range("D44").Value = "University"
$ws.Range("E44").Value = "Location"
$ws.Range("F44").Value = "image"

# --- Append a matching header row after the final blank separator row.
# After the insert above, the old last row (119, a blank "customFormat"
# spacer row) is now row 120, so the new header goes in row 121.
$ws.Range("B1:F1").Copy()
$ws.Range("B121:F121").PasteSpecial(-4122)

$ws.Range("B121").Value = "Name"
$ws.Range("C121").Value = "Credentials"
$ws.Range("D121").Value = "University"
$ws.Range("E121").Value = "Location"
$ws.Range("F121").Value = "image"

$excel.CutCopyMode = 0

# --- Update the active selection to match the post-edit view.
$ws.Range("B125").Select()
